$d = $word.ActiveDocument

$d.Content.Find.Execute("July 07, 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "July 30, 2020", 2)

$d.Content.Find.Execute("Pres. JOSE MARI L. CHAN", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MR. JONATHAN T. GOTIANUN", 2)

$d.Content.Find.Execute("Biscom Incorporated", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Davao Sugar Central Company, Inc.", 2)

$d.Content.Find.Execute("Unit 604, Legaspi Towers 200 Condominium, 107 Paseo de Roxas, Legaspi Vill., Makati City", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City", 2)

$d.Content.Find.Execute("Dear Pres. Chan:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dear Mr. Gotianun:", 2)

$d.Content.Find.Execute("Milling License No. 1001 for CY 2020 - 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Milling License No. 2021-01 for CY 2020 - 2021", 2)
